$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Convert E2/E3 from date-numbers to plain text strings matching the
# existing "dd-mm-yyyy" text style used elsewhere in column E. Forcing the
# cell to "Text" format before assigning the value keeps Excel from
# re-interpreting the string as a date; clearing formats afterwards drops
# the temporary text format so the cell is left with no explicit style,
# matching the other text cells in the column.
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "10-01-2029"
$ws.Range("E2").ClearFormats()

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "05-04-2030"
$ws.Range("E3").ClearFormats()

# Reset the active selection back to A1 (diff drops the stored <selection>
# pointing at G3).
$ws.Range("A1").Select()
